# Apply updated Cash Flow figures to the TWTR cash flow worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Change in payables and accrued liability
$ws.Range("B7").Value = 603382000.0
$ws.Range("C7").Value = 781000000.0
$ws.Range("D7").Value = 748000000.0
$ws.Range("E7").Value = 546097000.0
$ws.Range("F7").Value = 433375000.0
$ws.Range("G7").Value = 201378000.0

# Row 13: Long-Term Investments Change (Net)
$ws.Range("B13").Value = -54318000.0
$ws.Range("C13").Value = -13251000.0
$ws.Range("D13").Value = -56341000.0
$ws.Range("F13").Value = -53702000.0

# Row 18: Other financial activities
$ws.Range("B18").Value = 583378000.0

# Row 27: Investments Change (Net)
$ws.Range("B27").Value = -373166200.0
$ws.Range("F27").Value = -51558900.0

# Row 28: Issuance/Purchase of Shares
$ws.Range("B28").Value = -344278000.0

# Row 29: Capital Stock Change
$ws.Range("B29").Value = -344278000.0
